$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 1.53
$ws.Range("G2").Value = 1.65
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 9.6
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 5.1
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 3.35
$ws.Range("O2").Value = 1.34
$ws.Range("P2").Value = 1.77
$ws.Range("Q2").Value = 1.95
$ws.Range("R2").Value = 1.29
$ws.Range("S2").Value = 3.3
$ws.Range("T2").Value = 2.06
$ws.Range("U2").Value = 1.74
$ws.Range("V2").Value = 1.11
$ws.Range("W2").Value = 2.46
$ws.Range("H3").Value = 4.1
$ws.Range("I3").Value = 6.4
$ws.Range("J3").Value = 2.96
$ws.Range("K3").Value = 3.95
$ws.Range("Q3").Value = 2.18
$ws.Range("V3").Value = 1.19
$ws.Range("X3").Value = 14
$ws.Range("Z3").Value = 1000
$ws.Range("AA4").Value = 210
$ws.Range("AB4").Value = 8
$ws.Range("AC4").Value = 9.6
$ws.Range("AF4").Value = 12.5
$ws.Range("AG4").Value = 13
$ws.Range("AI4").Value = 140
$ws.Range("F4").Value = 1.78
$ws.Range("G4").Value = 1.92
$ws.Range("H4").Value = 5.1
$ws.Range("I4").Value = 6.2
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 3.8
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 2.54
$ws.Range("P4").Value = 1.6
$ws.Range("Q4").Value = 2.3
$ws.Range("R4").Value = 1.22
$ws.Range("S4").Value = 4.5
$ws.Range("T4").Value = 2.08
$ws.Range("U4").Value = 1.75
$ws.Range("W4").Value = 2.08
$ws.Range("X4").Value = 12.5
$ws.Range("AN5").Value = 120
$ws.Range("I5").Value = 1.54
$ws.Range("S5").Value = 2.98
$ws.Range("T5").Value = 1.94
$ws.Range("V5").Value = 2.86
$ws.Range("AI6").Value = 130
$ws.Range("AN6").Value = 3.3
$ws.Range("H6").Value = 13.5
$ws.Range("I6").Value = 14
$ws.Range("J6").Value = 7.4
$ws.Range("K6").Value = 7.6
$ws.Range("N6").Value = 8
$ws.Range("P6").Value = 3.4
$ws.Range("R6").Value = 1.94
$ws.Range("S6").Value = 2
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.07
$ws.Range("H7").Value = 2.2
$ws.Range("R7").Value = 1.5
$ws.Range("AA8").Value = 300
$ws.Range("AD8").Value = 32
$ws.Range("AE8").Value = 140
$ws.Range("AF8").Value = 8.199999999999999
$ws.Range("AI8").Value = 120
$ws.Range("AJ8").Value = 12.5
$ws.Range("AL8").Value = 38
$ws.Range("AC9").Value = 9.4
$ws.Range("M9").Value = 1.03
$ws.Range("AO10").Value = 17.5
$ws.Range("F10").Value = 2.28
$ws.Range("G10").Value = 2.32
$ws.Range("H10").Value = 3.15
$ws.Range("U10").Value = 2.8
$ws.Range("Y10").Value = 21
$ws.Range("F11").Value = 2.2
$ws.Range("G11").Value = 2.24
$ws.Range("S11").Value = 2.8
$ws.Range("T11").Value = 1.63
$ws.Range("AD12").Value = 70
$ws.Range("AF12").Value = 12.5
$ws.Range("AK12").Value = 12.5
$ws.Range("F12").Value = 1.17
$ws.Range("G12").Value = 1.18
$ws.Range("J12").Value = 10
$ws.Range("P12").Value = 4.4
$ws.Range("W12").Value = 6.4
$ws.Range("P13").Value = 2.34
$ws.Range("R13").Value = 1.54
$ws.Range("U13").Value = 2.52
